# Changes of 26th May 2022
# Update the ShipmentTracking (P), ActualRate (Q) and Result (R) columns on
# Sheet1 with the latest FedEx shipment tracking/result data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ShipmentTracking numbers for rows 2-26 (column P), in row order.
$trackingNumbers = @(
    "320018638745",
    "320018638756",
    "320018638789",
    "320018638804",
    "320018638848",
    "320018638860",
    "320018638892",
    "320018638918",
    "320018638940",
    "320018638962",
    "320018639009",
    "320018639020",
    "320018639053",
    "320018639075",
    "320018639101",
    "320018639123",
    "320018639167",
    "320018639189",
    "320018639215",
    "320018639237",
    "320018639260",
    "320018639270",
    "320018639281",
    "320018639292",
    "320018639307"
)

for ($i = 0; $i -lt $trackingNumbers.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Range("P$row")
    # Force the numeric-looking tracking number to be stored as text (the
    # column already holds text values), then drop the resulting
    # quote-prefix formatting so the cell keeps its original default style.
    $cell.Value2 = "'" + $trackingNumbers[$i]
    $cell.Style = "Normal"
}

# Row 24 now passes: the actual rate now matches the expected rate.
$q24 = $ws.Range("Q24")
$q24.Value2 = "'" + '$278.12'
$q24.Style = "Normal"

$ws.Range("R24").Value2 = "PASS"
